$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---- Latest Handback DateTime refresh ----
$zhcn.Range("L2").Value = "2017-02-22 08:44:49"
$zhcn.Range("L3").Value = "2017-02-22 08:44:49"

$dede.Range("L2").Value = "2017-02-22 08:45:13"
$dede.Range("L3").Value = "2017-02-22 08:45:13"

# ---- Error Detail cleared now that handback is in sync ----
$zhcn.Range("R3").Value = ""
$dede.Range("R3").Value = ""

# ---- Column width adjustments (status / error-detail columns) ----
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(18).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(18).ColumnWidth = 12.833333333333334

Write-Output "Generated handback report updates"
